$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update the li_partner id and the quick-placeholder date strings ---
$ws.Range("A2").Value = "z0bug.li_partner_6_2021"
$ws.Range("E2").Value = "<###-01-25"
$ws.Range("F2").Value = "<###-01-25"
$ws.Range("G2").Value = "<###-01-25"
$ws.Range("H2").Value = "<###-12-31"

# --- Row 3: new li_partner id for 2022 and corrected placeholder date ---
$ws.Range("A3").Value = "z0bug.li_partner_6_2022"
$ws.Range("E3").Value = "####-01-06"
$ws.Range("F3").Value = "####-01-06"
$ws.Range("G3").Value = "####-01-06"

# --- Sheet cosmetics: widen column A and move the active selection ---
$ws.Columns.Item(1).ColumnWidth = 21.15
$null = $ws.Range("A4").Select()
